$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking snapshot refresh: update Price (D) / Volume(1h) (E) for each
# row per the latest scrape, and fix the NEARProtocol/ImmutableX row order
# (rows 31-32 swapped places with refreshed figures).

$ws.Range("D2").Value = "68.927.23"
$ws.Range("E2").Value = "  +1.97%  "
$ws.Range("D3").Value = "3.735.75"
$ws.Range("E3").Value = "  -1.67%  "
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.15"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.67"
$ws.Range("E6").Value = "  -3.79%  "
$ws.Range("D7").Value = "3.731.82"
$ws.Range("E7").Value = "  -1.65%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +2.73%  "
$ws.Range("E10").Value = "  +3.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.34"
$ws.Range("E11").Value = "  +2.84%  "
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.06"
$ws.Range("E13").Value = "  +0.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000243"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").Value = "4.358.98"
$ws.Range("E15").Value = "  -1.90%  "
$ws.Range("D16").Value = "3.738.41"
$ws.Range("E16").Value = "  -2.70%  "
$ws.Range("D17").Value = "68.868.68"
$ws.Range("E17").Value = "  +1.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.25"
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.27"
$ws.Range("E20").Value = "  +6.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "497.76"
$ws.Range("E21").Value = "  +2.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.23"
$ws.Range("E22").Value = "  +12.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.722"
$ws.Range("E23").Value = "  -0.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.22"
$ws.Range("E24").Value = "  +1.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.30"
$ws.Range("E25").Value = "  -1.96%  "
$ws.Range("E26").Value = "  -6.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.36"
$ws.Range("E27").Value = "  +1.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.13"
$ws.Range("E28").Value = "  -0.86%  "
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.94"
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.44"
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.98"
$ws.Range("E32").Value = "  +3.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.75"
$ws.Range("E33").Value = "  -3.66%  "
$ws.Range("D34").Value = "3.885.65"
$ws.Range("E34").Value = "  -1.66%  "
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("D36").Value = "3.666.91"
$ws.Range("E36").Value = "  -1.97%  "
$ws.Range("E37").Value = "  -0.36%  "
$ws.Range("E38").Value = "  +1.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.81"
$ws.Range("E39").Value = "  +1.15%  "
$ws.Range("E40").Value = "  -1.87%  "
$ws.Range("E41").Value = "  -0.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "435.14"
$ws.Range("E42").Value = "  -4.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.91"
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("E44").Value = "  -0.51%  "
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("E46").Value = "  +1.30%  "
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.55"
$ws.Range("E48").Value = "  -2.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "142.20"
$ws.Range("E49").Value = "  +1.11%  "
$ws.Range("E50").Value = "  +0.92%  "
$ws.Range("D51").Value = "2.742.07"
$ws.Range("E51").Value = "  -2.99%  "
